$d = $word.ActiveDocument

# Locate the paragraph that marks the end of the bibliography section
# ("L.L. e Barros, M.T. org.. Drenagem Urbana. Ed. da Universidade e ABRH. 1995.").
# The three paragraphs that directly follow it:
#   1) an empty "Normal" paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. ..."
# must be removed entirely, leaving the bibliography paragraph directly followed
# by the existing blank paragraph that precedes the trailing page break.

$anchorText = "Drenagem Urbana. Ed. da Universidade e ABRH."
$removeFirstText = "Ver no Jupiter Salvar em pdf Salvar em docx"
$removeLastText = "Contact: luizeleno@usp.br"

$count = $d.Paragraphs.Count
$startPara = $null
$endPara = $null

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*$anchorText*") {
        # The paragraph right after this one starts the block to delete.
        $startPara = $i + 1
    }
    if ($t -like "*$removeFirstText*") {
        # Back up one paragraph to also catch the blank paragraph before it.
        $startPara = $i - 1
    }
    if ($t -like "*$removeLastText*") {
        $endPara = $i
    }
}

if ($startPara -ne $null -and $endPara -ne $null -and $endPara -ge $startPara) {
    $rangeStart = $d.Paragraphs.Item($startPara).Range.Start
    $rangeEnd = $d.Paragraphs.Item($endPara + 1).Range.Start
    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Delete()
    Write-Output "Deleted paragraphs $startPara through $endPara"
} else {
    Write-Output "Target paragraphs not found; no changes made."
}
